$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.310915589332581
$ws.Range("B1").Value = 2.08903980255127
$ws.Range("C1").Value = 5.093509197235107
$ws.Range("D1").Value = 1.985970139503479
$ws.Range("E1").Value = 1.079737544059753
